# Fruta / hortaliza, semanal
# The source row 19 (Comercializadora del Agro de Limari, Arandano (blue))
# is split into two weekly observations:
#   - row 19 is updated with the new week's data (date 45204 / 2023-10-05)
#   - row 20 is a new row that keeps the previous week's data that used to
#     live in row 19 (date 44455 / 2021-09-16)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, copy the existing row 19 down into the new row 20, preserving the
# date cell's number formatting.
$ws.Range("A20").Value = $ws.Range("A19").Value()
$ws.Range("B20").Value = $ws.Range("B19").Value()
$ws.Range("C20").Value = $ws.Range("C19").Value()

$ws.Range("D20").NumberFormat = $ws.Range("D19").NumberFormat
$ws.Range("D20").Value = $ws.Range("D19").Value()

$ws.Range("E20").Value = $ws.Range("E19").Value()
$ws.Range("F20").Value = $ws.Range("F19").Value()
$ws.Range("G20").Value = $ws.Range("G19").Value()
$ws.Range("H20").Value = $ws.Range("H19").Value()
$ws.Range("I20").Value = $ws.Range("I19").Value()
$ws.Range("J20").Value = $ws.Range("J19").Value()
$ws.Range("K20").Value = $ws.Range("K19").Value()
$ws.Range("L20").Value = $ws.Range("L19").Value()
$ws.Range("M20").Value = $ws.Range("M19").Value()
$ws.Range("N20").Value = $ws.Range("N19").Value()
$ws.Range("O20").Value = $ws.Range("O19").Value()
$ws.Range("P20").Value = $ws.Range("P19").Value()
$ws.Range("Q20").Value = $ws.Range("Q19").Value()
$ws.Range("R20").Value = $ws.Range("R19").Value()
$ws.Range("S20").Value = $ws.Range("S19").Value()
$ws.Range("T20").Value = $ws.Range("T19").Value()

# Now overwrite row 19 with the new week's data.
$ws.Range("A19").Value = 2
$ws.Range("B19").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C19").Value = "Coquimbo"
$ws.Range("D19").Value = 45204
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100101
$ws.Range("H19").Value = "Berries"
$ws.Range("I19").Value = 100101001
$ws.Range("J19").Value = "Arándano (blue)"
$ws.Range("K19").Value = "Sin especificar"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 400
$ws.Range("N19").Value = 9000
$ws.Range("O19").Value = 10000
$ws.Range("P19").Value = 9500
$ws.Range("Q19").Value = "$/bandeja 2 kilos"
$ws.Range("R19").Value = "Provincia de Limarí"
$ws.Range("S19").Value = 4750
$ws.Range("T19").Value = 2
